$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 170.33333
$ws.Range("J33").Value = 329.5
$ws.Range("L33").Value = 329.5
$ws.Range("N33").Value = -787.5

$ws.Range("H76").Value = 8709.652
$ws.Range("I76").Value = 13751.833
$ws.Range("J76").Value = 3209.0908
$ws.Range("K76").Value = 13751.833
$ws.Range("L76").Value = 3209.0908
$ws.Range("M76").Value = -13436.833
$ws.Range("N76").Value = -3839.0908

$ws.Range("H79").Value = 8709.652
$ws.Range("I79").Value = 13751.833
$ws.Range("J79").Value = 3209.0908
$ws.Range("K79").Value = 13751.833
$ws.Range("L79").Value = 3209.0908
$ws.Range("M79").Value = -12659.833
$ws.Range("N79").Value = -5393.0908

$ws.Range("H112").Value = 2353
$ws.Range("J112").Value = 2691.25
$ws.Range("L112").Value = 8073.75
$ws.Range("N112").Value = -10289.75

$ws.Range("H113").Value = 4349.143
$ws.Range("I113").Value = 2360
$ws.Range("J113").Value = 5454.222
$ws.Range("K113").Value = 2360
$ws.Range("L113").Value = 5454.222
$ws.Range("M113").Value = 894
$ws.Range("N113").Value = -11962.222

$ws.Range("H116").Value = 4582.0386
$ws.Range("I116").Value = 3363.3333
$ws.Range("J116").Value = 5626.643
$ws.Range("K116").Value = 3363.3333
$ws.Range("L116").Value = 5626.643
$ws.Range("M116").Value = 78.66670000000022
$ws.Range("N116").Value = -12510.643

$ws.Range("H132").Value = 2668.9062
$ws.Range("I132").Value = 1513.7587
$ws.Range("J132").Value = 13835.333
$ws.Range("K132").Value = 4541.2761
$ws.Range("L132").Value = 41505.999
$ws.Range("M132").Value = -2011.2761
$ws.Range("N132").Value = -46565.999

$ws.Range("H135").Value = 528790.0600000001
$ws.Range("I135").Value = 558156.2
$ws.Range("J135").Value = 200
$ws.Range("K135").Value = 5023405.8
$ws.Range("L135").Value = 1800
$ws.Range("M135").Value = -5020870.8
$ws.Range("N135").Value = -6870

$ws.Range("H138").Value = 2015.0676
$ws.Range("I138").Value = 1477.5186
$ws.Range("J138").Value = 2323.8723
$ws.Range("K138").Value = 4432.5558
$ws.Range("L138").Value = 6971.6169
$ws.Range("M138").Value = 707.4441999999999
$ws.Range("N138").Value = -17251.6169

$ws.Range("H141").Value = 2465.0833
$ws.Range("I141").Value = 2566.375
$ws.Range("J141").Value = 2262.5
$ws.Range("K141").Value = 7699.125
$ws.Range("L141").Value = 6787.5
$ws.Range("M141").Value = -2519.125
$ws.Range("N141").Value = -17147.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2083.3447
$ws.Range("I2").Value = 1310.381
$ws.Range("J2").Value = 4112.375
$ws.Range("K2").Value = 1310.381
$ws.Range("L2").Value = 4112.375
$ws.Range("M2").Value = -1197.381
$ws.Range("N2").Value = -4338.375

$ws.Range("H45").Value = 3899
$ws.Range("I45").Value = 3110.1538
$ws.Range("J45").Value = 5950
$ws.Range("K45").Value = 3110.1538
$ws.Range("L45").Value = 5950
$ws.Range("M45").Value = -2733.1538
$ws.Range("N45").Value = -6704

$ws.Range("H61").Value = 1764.5927
$ws.Range("I61").Value = 1727
$ws.Range("K61").Value = 1727
$ws.Range("M61").Value = -1515

$ws.Range("H116").Value = 2083.3447
$ws.Range("I116").Value = 1310.381
$ws.Range("J116").Value = 4112.375
$ws.Range("K116").Value = 1310.381
$ws.Range("L116").Value = 4112.375
$ws.Range("M116").Value = 983.6189999999999
$ws.Range("N116").Value = -8700.375

$ws.Range("H136").Value = 1764.5927
$ws.Range("I136").Value = 1727
$ws.Range("K136").Value = 5181
$ws.Range("M136").Value = -2631

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2083.3447
$ws.Range("I3").Value = 1310.381
$ws.Range("J3").Value = 4112.375
$ws.Range("K3").Value = 1310.381
$ws.Range("L3").Value = 4112.375
$ws.Range("M3").Value = -1196.381
$ws.Range("N3").Value = -4340.375

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5125
$ws.Range("I16").Value = 5125
$ws.Range("K16").Value = 5125
$ws.Range("M16").Value = -4838

$ws.Range("H31").Value = 3365.3022
$ws.Range("I31").Value = 2244.5
$ws.Range("J31").Value = 4029.4814
$ws.Range("K31").Value = 2244.5
$ws.Range("L31").Value = 4029.4814
$ws.Range("M31").Value = -1949.5
$ws.Range("N31").Value = -4619.481400000001

$ws.Range("H34").Value = 3365.3022
$ws.Range("I34").Value = 2244.5
$ws.Range("J34").Value = 4029.4814
$ws.Range("K34").Value = 2244.5
$ws.Range("L34").Value = 4029.4814
$ws.Range("M34").Value = -2042.5
$ws.Range("N34").Value = -4433.481400000001

$ws.Range("H41").Value = 17446.666
$ws.Range("J41").Value = 19956
$ws.Range("L41").Value = 19956
$ws.Range("N41").Value = -20812

$ws.Range("H48").Value = 5993.3335
$ws.Range("J48").Value = 5993.3335
$ws.Range("L48").Value = 5993.3335
$ws.Range("N48").Value = -6945.3335

$ws.Range("H50").Value = 9086
$ws.Range("J50").Value = 9086
$ws.Range("L50").Value = 9086
$ws.Range("N50").Value = -10336

$ws.Range("H58").Value = 1540.6666
$ws.Range("I58").Value = 1834.5625
$ws.Range("J58").Value = 952.875
$ws.Range("K58").Value = 1834.5625
$ws.Range("L58").Value = 952.875
$ws.Range("M58").Value = -1631.5625
$ws.Range("N58").Value = -1358.875

$ws.Range("H60").Value = 30516
$ws.Range("J60").Value = 30516
$ws.Range("L60").Value = 30516
$ws.Range("N60").Value = -31538

$ws.Range("H109").Value = 10520
$ws.Range("J109").Value = 10520
$ws.Range("L109").Value = 10520
$ws.Range("N109").Value = -12600

$ws.Range("H113").Value = 5125
$ws.Range("I113").Value = 5125
$ws.Range("K113").Value = 5125
$ws.Range("M113").Value = -2955

$ws.Range("H136").Value = 1540.6666
$ws.Range("I136").Value = 1834.5625
$ws.Range("J136").Value = 952.875
$ws.Range("K136").Value = 5503.6875
$ws.Range("L136").Value = 2858.625
$ws.Range("M136").Value = -2953.6875
$ws.Range("N136").Value = -7958.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 365
$ws.Range("I5").Value = 347.91666
$ws.Range("K5").Value = 1043.74998
$ws.Range("M5").Value = -931.7499800000001

$ws.Range("H38").Value = 120.3
$ws.Range("I38").Value = 58.5
$ws.Range("J38").Value = 213
$ws.Range("K38").Value = 175.5
$ws.Range("L38").Value = 639
$ws.Range("M38").Value = 171.5
$ws.Range("N38").Value = -1333

$ws.Range("H132").Value = 824450.1
$ws.Range("I132").Value = 1013030.94
$ws.Range("K132").Value = 9117278.459999999
$ws.Range("M132").Value = -9114748.459999999

$ws.Range("H135").Value = 365
$ws.Range("I135").Value = 347.91666
$ws.Range("K135").Value = 3131.24994
$ws.Range("M135").Value = -596.2499399999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 5266.6665
$ws.Range("I19").Value = 6066.6665
$ws.Range("K19").Value = 6066.6665
$ws.Range("M19").Value = -5778.6665

$ws.Range("H46").Value = 10482.667
$ws.Range("J46").Value = 11728.333
$ws.Range("L46").Value = 11728.333
$ws.Range("N46").Value = -12040.333

$ws.Range("H80").Value = 3248
$ws.Range("I80").Value = 2896
$ws.Range("J80").Value = 3600
$ws.Range("K80").Value = 2896
$ws.Range("L80").Value = 3600
$ws.Range("M80").Value = -1898
$ws.Range("N80").Value = -5596

$ws.Range("H83").Value = 3248
$ws.Range("I83").Value = 2896
$ws.Range("J83").Value = 3600
$ws.Range("K83").Value = 14480
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -9488
$ws.Range("N83").Value = -27984

$ws.Range("H123").Value = 34309.5
$ws.Range("J123").Value = 34309.5
$ws.Range("L123").Value = 34309.5
$ws.Range("N123").Value = -39209.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 36150
$ws.Range("J64").Value = 36150
$ws.Range("L64").Value = 36150
$ws.Range("N64").Value = -36600

$ws.Range("H67").Value = 36150
$ws.Range("J67").Value = 36150
$ws.Range("L67").Value = 36150
$ws.Range("N67").Value = -37710

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 19700
$ws.Range("J109").Value = 19700
$ws.Range("L109").Value = 19700
$ws.Range("N109").Value = -22474

$ws.Range("H132").Value = 3591.3635
$ws.Range("I132").Value = 4320.143
$ws.Range("J132").Value = 2316
$ws.Range("K132").Value = 12960.429
$ws.Range("L132").Value = 6948
$ws.Range("M132").Value = -12008
